$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.414.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '''2.092.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''329.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").Value = '''0.5205'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''0.4335'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").Value = '''52.07'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +15.39%  '
$ws.Range("D10").Value = '''0.08831'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").Value = '''1.155'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("D12").Value = '''24.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").Value = '''2.087.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").Value = '''6.674'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.80%  '
$ws.Range("D15").Value = '''7.666'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '''95.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").Value = '''1.004'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").Value = '''0.00001117'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.19%  '
$ws.Range("D19").Value = '''0.06583'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").Value = '''19.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").Value = '''1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '''6.263'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").Value = '''30.471.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.54%  '
$ws.Range("D24").Value = '''12.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").Value = '''2.339'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.73%  '
$ws.Range("D26").Value = '''2.332.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.50%  '
$ws.Range("D27").Value = '''22.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = '''2.584'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").Value = '''162.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '''131.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("D31").Value = '''1.188'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = '''0.1067'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").Value = '''1.665'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.06%  '
$ws.Range("D34").Value = '''6.133'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.78%  '
$ws.Range("D35").Value = '''3.885'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").Value = '''10.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.51%  '
$ws.Range("D37").Value = '''0.02562'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").Value = '''0.06791'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '''5.454'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").Value = '''12.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("D41").Value = '''0.2261'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '''0.6890'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").Value = '''1.265'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").Value = '''1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '''0.6364'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").Value = '''13.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.79%  '
$ws.Range("D47").Value = '''2.197'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").Value = '''3.621'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("D49").Value = '''1.233'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.49%  '
$ws.Range("D50").Value = '''1.239'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("D51").Value = '''81.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.90%  '
